$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: the "Nb nouveaux deces a l'hopital" (L) and "...extra-hospitaliers" (M)
# columns are formatted as Text (@), so a plain .Value write of a number gets
# stored as text. Flip the format to General for the write, then restore it,
# so the stored cell keeps its original Text display format but a genuine
# numeric value (matching how these columns already look further up the sheet).
function Set-NumericValue($range, $value) {
    $range.NumberFormat = "General"
    $range.Value = $value
    $range.NumberFormat = "@"
}

# Revised "new positive cases" counts for a few already-recorded days
$ws.Range("C350").Value = 56
$ws.Range("C351").Value = 61
$ws.Range("C352").Value = 59
Set-NumericValue $ws.Range("L352") 2

# Fill in row 353 (2021-02-05, serial 44239) with that day's figures
$ws.Range("C353").Value = 12
$ws.Range("E353").Value = 11
$ws.Range("F353").Value = 7
$ws.Range("G353").Value = 74
Set-NumericValue $ws.Range("L353") 1
Set-NumericValue $ws.Range("M353") 0

# Reset the saved selection on the frozen pane to A2
$ws.Range("A2").Select()
